# Update "想去人数" (want-to-go count) figures refreshed at the
# gh-pages data snapshot (commit 7921097).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 996
$ws.Range("F5").Value  = 2
$ws.Range("F6").Value  = 7188
$ws.Range("F8").Value  = 921
$ws.Range("F10").Value = 764
$ws.Range("F11").Value = 546
$ws.Range("F12").Value = 64
$ws.Range("F15").Value = 836
$ws.Range("F16").Value = 2868
$ws.Range("F17").Value = 149
$ws.Range("F18").Value = 33
$ws.Range("F19").Value = 229
$ws.Range("F20").Value = 736
$ws.Range("F22").Value = 426
$ws.Range("F23").Value = 17
$ws.Range("F24").Value = 139
$ws.Range("F25").Value = 185
$ws.Range("F26").Value = 129
$ws.Range("F27").Value = 187
$ws.Range("F29").Value = 67
$ws.Range("F30").Value = 180
$ws.Range("F33").Value = 293
$ws.Range("F34").Value = 374
$ws.Range("F38").Value = 37

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 194

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 194
$ws.Range("F3").Value  = 996
$ws.Range("F9").Value  = 2
$ws.Range("F10").Value = 7188
$ws.Range("F12").Value = 921
$ws.Range("F14").Value = 764
$ws.Range("F15").Value = 546
$ws.Range("F16").Value = 64
$ws.Range("F19").Value = 836
$ws.Range("F21").Value = 2868
$ws.Range("F22").Value = 149
$ws.Range("F23").Value = 33
$ws.Range("F25").Value = 229
$ws.Range("F26").Value = 736
$ws.Range("F29").Value = 426
$ws.Range("F30").Value = 17
$ws.Range("F31").Value = 139
$ws.Range("F32").Value = 185
$ws.Range("F33").Value = 129
$ws.Range("F34").Value = 187
$ws.Range("F36").Value = 67
$ws.Range("F37").Value = 180
$ws.Range("F40").Value = 293
$ws.Range("F41").Value = 374
$ws.Range("F45").Value = 37
